$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1224.3572
$ws.Range("I43").Value = 1000.5
$ws.Range("J43").Value = 1261.6666
$ws.Range("K43").Value = 1000.5
$ws.Range("L43").Value = 1261.6666
$ws.Range("M43").Value = -931.5
$ws.Range("N43").Value = -1399.6666

$ws.Range("H125").Value = 3833.5
$ws.Range("I125").Value = 4681
$ws.Range("J125").Value = 3197.875
$ws.Range("K125").Value = 42129
$ws.Range("L125").Value = 28780.875
$ws.Range("M125").Value = -39669
$ws.Range("N125").Value = -33700.875

$ws.Range("H132").Value = 909.9643
$ws.Range("I132").Value = 844.73914
$ws.Range("J132").Value = 1210
$ws.Range("K132").Value = 2534.21742
$ws.Range("L132").Value = 3630
$ws.Range("M132").Value = -4.217419999999947
$ws.Range("N132").Value = -8690

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3956.5676
$ws.Range("I32").Value = 3048.2295
$ws.Range("J32").Value = 8218.77
$ws.Range("K32").Value = 3048.2295
$ws.Range("L32").Value = 8218.77
$ws.Range("M32").Value = -2761.2295
$ws.Range("N32").Value = -8792.77

$ws.Range("H61").Value = 5894.087
$ws.Range("I61").Value = 6071.091
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 6071.091
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -5859.091
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 1060.7028
$ws.Range("I74").Value = 1172.76
$ws.Range("J74").Value = 827.25
$ws.Range("K74").Value = 1172.76
$ws.Range("L74").Value = 827.25
$ws.Range("M74").Value = -298.76
$ws.Range("N74").Value = -2575.25

$ws.Range("H77").Value = 1060.7028
$ws.Range("I77").Value = 1172.76
$ws.Range("J77").Value = 827.25
$ws.Range("K77").Value = 5863.8
$ws.Range("L77").Value = 4136.25
$ws.Range("M77").Value = -1495.8
$ws.Range("N77").Value = -12872.25

$ws.Range("H97").Value = 740.0909
$ws.Range("I97").Value = 613
$ws.Range("J97").Value = 2011
$ws.Range("K97").Value = 613
$ws.Range("L97").Value = 2011
$ws.Range("M97").Value = -117
$ws.Range("N97").Value = -3003

$ws.Range("H122").Value = 1426080.6
$ws.Range("I122").Value = 1710429.1
$ws.Range("J122").Value = 4338
$ws.Range("K122").Value = 5131287.300000001
$ws.Range("L122").Value = 13014
$ws.Range("M122").Value = -5128837.300000001
$ws.Range("N122").Value = -17914

$ws.Range("H132").Value = 3098.6743
$ws.Range("I132").Value = 1766.6897
$ws.Range("K132").Value = 5300.0691
$ws.Range("M132").Value = -2770.0691

$ws.Range("H136").Value = 5894.087
$ws.Range("I136").Value = 6071.091
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 18213.273
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -15663.273
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 36435.75
$ws.Range("J62").Value = 36435.75
$ws.Range("L62").Value = 36435.75
$ws.Range("N62").Value = -37807.75

$ws.Range("H65").Value = 36435.75
$ws.Range("J65").Value = 36435.75
$ws.Range("L65").Value = 109307.25
$ws.Range("N65").Value = -116171.25

$ws.Range("H94").Value = 1055.1143
$ws.Range("I94").Value = 626.0357
$ws.Range("K94").Value = 626.0357
$ws.Range("M94").Value = -175.0357

$ws.Range("H105").Value = 2435.75
$ws.Range("I105").Value = 1691.6154
$ws.Range("J105").Value = 3817.7144
$ws.Range("K105").Value = 1691.6154
$ws.Range("L105").Value = 3817.7144
$ws.Range("M105").Value = 55.38460000000009
$ws.Range("N105").Value = -7311.7144

$ws.Range("H107").Value = 1296.2858
$ws.Range("I107").Value = 1432.2
$ws.Range("J107").Value = 956.5
$ws.Range("K107").Value = 1432.2
$ws.Range("L107").Value = 956.5
$ws.Range("M107").Value = 487.8
$ws.Range("N107").Value = -4796.5

$ws.Range("H134").Value = 4374.15
$ws.Range("I134").Value = 4723.9375
$ws.Range("J134").Value = 2975
$ws.Range("K134").Value = 14171.8125
$ws.Range("L134").Value = 8925
$ws.Range("M134").Value = -11636.8125
$ws.Range("N134").Value = -13995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 193.45454
$ws.Range("I7").Value = 196
$ws.Range("K7").Value = 196
$ws.Range("M7").Value = -83

$ws.Range("H31").Value = 3047.2273
$ws.Range("I31").Value = 1512.4762
$ws.Range("J31").Value = 4448.522
$ws.Range("K31").Value = 1512.4762
$ws.Range("L31").Value = 4448.522
$ws.Range("M31").Value = -1217.4762
$ws.Range("N31").Value = -5038.522

$ws.Range("H34").Value = 3047.2273
$ws.Range("I34").Value = 1512.4762
$ws.Range("J34").Value = 4448.522
$ws.Range("K34").Value = 1512.4762
$ws.Range("L34").Value = 4448.522
$ws.Range("M34").Value = -1310.4762
$ws.Range("N34").Value = -4852.522

$ws.Range("H58").Value = 1392.6842
$ws.Range("I58").Value = 1037.1428
$ws.Range("J58").Value = 1831.8823
$ws.Range("K58").Value = 1037.1428
$ws.Range("L58").Value = 1831.8823
$ws.Range("M58").Value = -834.1428000000001
$ws.Range("N58").Value = -2237.8823

$ws.Range("H99").Value = 2972.258
$ws.Range("I99").Value = 2035.1538
$ws.Range("J99").Value = 7845.2
$ws.Range("K99").Value = 2035.1538
$ws.Range("L99").Value = 7845.2
$ws.Range("M99").Value = -537.1538
$ws.Range("N99").Value = -10841.2

$ws.Range("H126").Value = 2972.258
$ws.Range("I126").Value = 2035.1538
$ws.Range("J126").Value = 7845.2
$ws.Range("K126").Value = 6105.4614
$ws.Range("L126").Value = 23535.6
$ws.Range("M126").Value = -3635.4614
$ws.Range("N126").Value = -28475.6

$ws.Range("H136").Value = 1392.6842
$ws.Range("I136").Value = 1037.1428
$ws.Range("J136").Value = 1831.8823
$ws.Range("K136").Value = 3111.4284
$ws.Range("L136").Value = 5495.6469
$ws.Range("M136").Value = -561.4284000000002
$ws.Range("N136").Value = -10595.6469

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 231937.84
$ws.Range("I5").Value = 537.2308
$ws.Range("J5").Value = 463338.47
$ws.Range("K5").Value = 1611.6924
$ws.Range("L5").Value = 1390015.41
$ws.Range("M5").Value = -1499.6924
$ws.Range("N5").Value = -1390239.41

$ws.Range("H131").Value = 1064781.4
$ws.Range("I131").Value = 5882963
$ws.Range("J131").Value = 1027.1039
$ws.Range("K131").Value = 17648889
$ws.Range("L131").Value = 3081.3117
$ws.Range("M131").Value = -17643849
$ws.Range("N131").Value = -13161.3117

$ws.Range("H132").Value = 2176.125
$ws.Range("J132").Value = 2601
$ws.Range("L132").Value = 23409
$ws.Range("N132").Value = -28469

$ws.Range("H135").Value = 231937.84
$ws.Range("I135").Value = 537.2308
$ws.Range("J135").Value = 463338.47
$ws.Range("K135").Value = 4835.077200000001
$ws.Range("L135").Value = 4170046.23
$ws.Range("M135").Value = -2300.077200000001
$ws.Range("N135").Value = -4175116.23

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 817
$ws.Range("I97").Value = 778.8889
$ws.Range("J97").Value = 902.75
$ws.Range("K97").Value = 778.8889
$ws.Range("L97").Value = 902.75
$ws.Range("M97").Value = -282.8889
$ws.Range("N97").Value = -1894.75

$ws.Range("H102").Value = 1351
$ws.Range("I102").Value = 1222.4
$ws.Range("J102").Value = 1479.6
$ws.Range("K102").Value = 1222.4
$ws.Range("L102").Value = 1479.6
$ws.Range("M102").Value = 399.5999999999999
$ws.Range("N102").Value = -4723.6

$ws.Range("H126").Value = 6666.524
$ws.Range("I126").Value = 9367.23
$ws.Range("J126").Value = 2277.875
$ws.Range("K126").Value = 28101.69
$ws.Range("L126").Value = 6833.625
$ws.Range("M126").Value = -25631.69
$ws.Range("N126").Value = -11773.625

$ws.Range("H132").Value = 2880.4082
$ws.Range("I132").Value = 3445.7273
$ws.Range("J132").Value = 2419.7778
$ws.Range("K132").Value = 10337.1819
$ws.Range("L132").Value = 7259.3334
$ws.Range("M132").Value = -7807.1819
$ws.Range("N132").Value = -12319.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 27789366
$ws.Range("I93").Value = 17216.334
$ws.Range("J93").Value = 83333660
$ws.Range("K93").Value = 17216.334
$ws.Range("L93").Value = 83333660
$ws.Range("M93").Value = -15968.334
$ws.Range("N93").Value = -83336156

$ws.Range("H122").Value = 6265351.5
$ws.Range("I122").Value = 7144557
$ws.Range("J122").Value = 3334666.8
$ws.Range("K122").Value = 21433671
$ws.Range("L122").Value = 10004000.4
$ws.Range("M122").Value = -21431221
$ws.Range("N122").Value = -10008900.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1578
$ws.Range("I122").Value = 1863.3334
$ws.Range("J122").Value = 1150
$ws.Range("K122").Value = 5590.0002
$ws.Range("L122").Value = 3450
$ws.Range("M122").Value = -3140.0002
$ws.Range("N122").Value = -8350

$ws.Range("H132").Value = 1877.6
$ws.Range("I132").Value = 1137.1428
$ws.Range("K132").Value = 3411.4284
$ws.Range("M132").Value = -881.4284000000002
